# fix(publipostage): Correct status name
#
# The "statut_label" column (B) uses a short color-code word and the
# "statut_name" column (C) uses a longer human-readable description.
# This corrects the wording of one color label and of the four
# result/publication status descriptions, wherever they occur.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "bleu" = "noir";
    "résultat et / ou publication posté" = "résultat postés ou publiés";
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés";
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois";
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $current = $cell.Value2
        if ($replacements.ContainsKey($current)) {
            $cell.Value = $replacements[$current]
        }
    }
}
